$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    $ws.Cells.Item($r, 4).Value = 0.018                         # D - historical_growth_revenue_last_5_years
    $ws.Cells.Item($r, 5).Value = -0.104                        # E - historical_growth_net_income_last_5_years (new)
    $ws.Cells.Item($r, 7).Value = 0.1115220483641536            # G - ebitdard_margin
    $ws.Cells.Item($r, 8).Value = 0.1115220483641536            # H - ebitda_margin
    $ws.Cells.Item($r, 9).Value = 0.1243243243243243            # I - operating_margin
    $ws.Cells.Item($r, 10).Value = 0.1119174730285842           # J - after_tax_operating_margin
    $ws.Cells.Item($r, 11).Value = 1.38                         # K - trailing_net_income
    $ws.Cells.Item($r, 12).Value = 0.01963015647226174          # L - net_margin
    $ws.Cells.Item($r, 21).Value = 3.86                         # U - cash
    $ws.Cells.Item($r, 22).Value = 0.03580705009276438          # V - cash_market_cap
    $ws.Cells.Item($r, 23).Value = 0.02827868852459016          # W - roe
    $ws.Cells.Item($r, 24).Value = 0.07579399897851161          # X - cost_equity
    $ws.Cells.Item($r, 25).Value = -0.04751531045392145         # Y - roe_cost_equity
    $ws.Cells.Item($r, 26).Value = 1.608695652173913            # Z - sales_invested_capital
    $ws.Cells.Item($r, 27).Value = 0.1800411522633745           # AA - roic
    $ws.Cells.Item($r, 28).Value = 0.07556846554090965          # AB - cost_capital
    $ws.Cells.Item($r, 29).Value = 0.1044726867224649           # AC - roic_cost_capital
    $ws.Cells.Item($r, 30).Value = 0.707                        # AD - debt_total
    $ws.Cells.Item($r, 31).Value = 0                            # AE - debt_leases
    $ws.Cells.Item($r, 32).Value = 0.707                        # AF - debt_total_inc_leases
    $ws.Cells.Item($r, 33).Value = -3.153                       # AG - net_debt
    $ws.Cells.Item($r, 34).Value = 0.006515708663957165         # AH - debt_market_capital
    $ws.Cells.Item($r, 35).Value = 0.01606562592314859          # AI - debt_book_capital
    $ws.Cells.Item($r, 36).Value = -0.03012986516574771         # AJ - net_debt_market_capital
    $ws.Cells.Item($r, 37).Value = -0.07853637880788104         # AK - net_debt_book_capital
    $ws.Cells.Item($r, 38).Value = 0.081                        # AL - interest_expenses
    $ws.Cells.Item($r, 39).Value = 0.081                        # AM - net_interest_expenses
    $ws.Cells.Item($r, 40).Value = 0.07403141361256543          # AN - debt_ebitda
    $ws.Cells.Item($r, 41).Value = 107.9012345679012            # AO - ebit_interest_expenses
    $ws.Cells.Item($r, 42).Value = -0.3301570680628272          # AP - net_debt_ebitda
    $ws.Cells.Item($r, 43).Value = 107.9012345679012            # AQ - ebit_net_interest_expenses
}

# Row 2: M..S all become 0 (positive zero)
$ws.Cells.Item(2, 13).Value = 0   # M
$ws.Cells.Item(2, 14).Value = 0   # N
$ws.Cells.Item(2, 15).Value = 0   # O
$ws.Cells.Item(2, 16).Value = 0   # P
$ws.Cells.Item(2, 17).Value = 0   # Q
$ws.Cells.Item(2, 18).Value = 0   # R
$ws.Cells.Item(2, 19).Value = 0   # S

# Row 3: M..R become -0 (negative zero), S becomes 0 (positive zero)
$ws.Cells.Item(3, 13).Value = -0   # M
$ws.Cells.Item(3, 14).Value = -0   # N
$ws.Cells.Item(3, 15).Value = -0   # O
$ws.Cells.Item(3, 16).Value = -0   # P
$ws.Cells.Item(3, 17).Value = -0   # Q
$ws.Cells.Item(3, 18).Value = -0   # R
$ws.Cells.Item(3, 19).Value = 0    # S

# T column (buybacks_cash_returned) removed entirely for rows 2 and 3
$ws.Cells.Item(2, 20).ClearContents()
$ws.Cells.Item(3, 20).ClearContents()
